$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - sample sizes changed
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 - CON data tweaks
$ws.Range("B2").Value = 42.225061315727089
$ws.Range("C2").Value = 39.71845511051913
$ws.Range("D2").Value = 46.776653437843954
$ws.Range("E2").Value = 29.200924079334971

# Row 3 - STR data tweaks
$ws.Range("B3").Value = 57.342222431918422
$ws.Range("C3").Value = 38.427319165220226
$ws.Range("D3").Value = 42.86379260784507
$ws.Range("E3").Value = 17.26090334741032

# Update the selection to match the narrower highlighted range
$ws.Range("B1:E3").Select()
